# resumeTemplate.docx: collapse the old "Work Experience" block (job
# history lines + laugh/live/love/eat/sleep/be-merry bullets +
# "Extracurriculars" block + never/ever/cry bullets) down to just a
# "Work Experience:" heading followed by a two-run "Shawmut"/"Workman"
# line, per "Store items as JSON in resumeWrite.py".

$d = $word.ActiveDocument

# 1) Drop every paragraph after the first one (the bullet lists and the
#    Extracurriculars paragraph all disappear).
if ($d.Paragraphs.Count -gt 1) {
    $tail = $d.Range($d.Paragraphs(2).Range.Start, $d.Paragraphs($d.Paragraphs.Count).Range.End)
    $tail.Delete()
}

# 2) Clear out the remaining paragraph's runs (keep the paragraph mark)
#    so we can rebuild it without inheriting the old bold "Work
#    Experience:" run formatting.
$p1 = $d.Paragraphs(1)
$p1Body = $d.Range($p1.Range.Start, $p1.Range.End - 1)
if ($p1Body.Start -ne $p1Body.End) {
    $p1Body.Delete()
}

# 3) Paragraph 1 becomes a Heading 2 "Work Experience:" line (plain run;
#    the heading style itself carries the bold weight).
$p1.Style = "Heading 2"
$p1.Range.InsertAfter("Work Experience:")

# 4) New second paragraph, back to the (implicit) Normal body style.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Style = "Normal"

# 5) "Shawmut" (bold) followed by "Workman" (plain) as two runs.
$p2.Range.InsertAfter("ShawmutWorkman")
$boldRange = $d.Range($p2.Range.Start, $p2.Range.Start + 7)
$boldRange.Font.Bold = 1
